# df_RSE_settings.xlsx edit
# Commit: "Created a copy of the RSE settings file (i.e. ORIGINAL)."
#
# Semantic change: the "Minimum" (col B) / "Maximum" (col C) coefficient
# bounds for every material row (rows 11-37) are reset to a uniform
# Min=1 / Max=0.25, replacing the previous per-row values (2/0.5 or 5/1).
# The active selection is also reset back to the top of the sheet (B1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 11 through 37: column B -> 1, column C -> 0.25
$ws.Range("B11:B37").Value = 1
$ws.Range("C11:C37").Value = 0.25

# Reset the selection/active cell to B1 (was B9)
$ws.Range("B1").Select() | Out-Null
